$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "Strong use of calculus (differentiation/chain rule), series (AP/GP sums), integration with definite integrals and limits, and coordinate geometry. Trigonometric equation solving to final angles was correct. Good handling of inverse functions and function range concepts."
$ws.Range("L2").Value = "Transformation descriptions sometimes lacked clear order; occasional algebraic slips (e.g. Q8(b) perimeter leading to incorrect r). In Q9(a) sign/interval reasoning led to an incorrect final interval. In Q10(b) an incorrect y-value for point B propagated through the perpendicular bisector; final equation not simplified to ax+by+c=0. Careful checking of values and clarity of statements would improve accuracy."

$ws.Range("M2").Value = 64

$ws.Range("C3").Value = 4
$ws.Range("C5").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 3
